$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of column letter -> new value, applied identically to rows 2 and 3
$updates = @{
    "D"  = -0.61
    "G"  = -146.7741935483871
    "H"  = -429.0322580645162
    "I"  = -395.9763829178406
    "J"  = -395.9763829178406
    "K"  = -12.4
    "L"  = -400
    "U"  = 5.37
    "V"  = 0.02799791449426486
    "W"  = -2.194690265486726
    "X"  = 0.08328667581471393
    "Y"  = -2.277976941301439
    "Z"  = -0.02744186943412057
    "AA" = 10.86633219902671
    "AB" = 0.08324827837125108
    "AC" = 10.78308392065546
    "AD" = 0.092
    "AE" = 0.1463393522652885
    "AF" = 0.2383393522652885
    "AG" = -5.131660647734711
    "AH" = 0.001241102964487164
    "AI" = 0.0584402943646424
    "AJ" = -0.02749079284436478
    "AK" = 3.972917079059825
    "AL" = 0.044
    "AM" = -0.02600000000000001
    "AN" = -0.007217950729640672
    "AO" = -281.8181818181819
    "AP" = 0.4026094969194031
    "AQ" = 476.9230769230768
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("$col" + "2").Value = $value
    $ws.Range("$col" + "3").Value = $value
}
